# Rename the worksheet tab (was "alpha3F-HW50.xpc" -> "alpha3F")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "alpha3F"

# Append a new data row (row 16) that mirrors the layout/formatting of the
# preceding row (row 15), as part of exporting the new Gaussian Quadrature
# Scheme data into the averaged-intensities worksheet.
$ws.Range("A15:M15").Copy()
$ws.Range("A16:M16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:M16").Value = 1
